# =====================================================================
# Commit: "output of tailing.xlsx output of intensities.xlsx using numba
# to speed up age calculation added specific constants to constants
# sheet in Results.xlsx added support for string labnrs"
#
# This script rewrites:
#   - Constants sheet: reorders/relabels all existing rows and appends
#     11 new constant rows (standard weights, blank/tail-correction
#     constants, a "type" field, etc.)
#   - Calc + Results sheets: updated numeric outputs for rows 7/9/11/13
#     (re-run of the age-calculation pipeline)
# =====================================================================

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) Constants sheet
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Constants")

# Extend formatting (style) down to the 11 new rows (39-49) by copying
# the format of the last pre-existing row (38) down the column.
$ws.Range("A38:B38").Copy()
$ws.Range("A39:B49").PasteSpecial(-4122)

# Column A is slightly narrower now.
$ws.Columns.Item(1).ColumnWidth = 16.83

# Row 1: Blank
$ws.Range("A1").Value = "Blank"
$ws.Range("B1").Value = $true
# Row 2: Yield_U
$ws.Range("A2").Value = "Yield_U"
$ws.Range("B2").Value = 1
# Row 3: Yield_Th
$ws.Range("A3").Value = "Yield_Th"
$ws.Range("B3").Value = 1
# Row 4: Gain
$ws.Range("A4").Value = "Gain"
$ws.Range("B4").Value = 1
# Row 5: Tail shift
$ws.Range("A5").Value = "Tail shift"
$ws.Range("B5").Value = $false
# Row 6: mf48
$ws.Range("A6").Value = "mf48"
$ws.Range("B6").Value = 1.336402435064349
# Row 7: mf36
$ws.Range("A7").Value = "mf36"
$ws.Range("B7").Value = 1.008202776684838
# Row 8: mf56
$ws.Range("A8").Value = "mf56"
$ws.Range("B8").Value = 0.334493224630051
# Row 9: mf68
$ws.Range("A9").Value = "mf68"
$ws.Range("B9").Value = 0.665506775369946
# Row 10: mf92
$ws.Range("A10").Value = "mf92"
$ws.Range("B10").Value = 1.025840620457897
# Row 11: mf38
$ws.Range("A11").Value = "mf38"
$ws.Range("B11").Value = 1.673784240557133
# Row 12: mf35
$ws.Range("A12").Value = "mf35"
$ws.Range("B12").Value = 0.673784240557127
# Row 13: mf43
$ws.Range("A13").Value = "mf43"
$ws.Range("B13").Value = -0.337307116990441
# Row 14: mf45
$ws.Range("A14").Value = "mf45"
$ws.Range("B14").Value = 0.336402435064353
# Row 15: mf09
$ws.Range("A15").Value = "mf09"
$ws.Range("B15").Value = -0.34318587041139
# Row 16: mf29
$ws.Range("A16").Value = "mf29"
$ws.Range("B16").Value = -1.025840620457897
# Row 17: mf34
$ws.Range("A17").Value = "mf34"
$ws.Range("B17").Value = 0.337307116990439
# Row 18: mf58
$ws.Range("A18").Value = "mf58"
$ws.Range("B18").Value = 1
# Row 19: mf02
$ws.Range("A19").Value = "mf02"
$ws.Range("B19").Value = 0.682654750046506
# Row 20: l230
$ws.Range("A20").Value = "l230"
$ws.Range("B20").Value = [double]"9.1705E-06"
# Row 21: l232
$ws.Range("A21").Value = "l232"
$ws.Range("B21").Value = [double]"4.94752E-11"
# Row 22: l234
$ws.Range("A22").Value = "l234"
$ws.Range("B22").Value = [double]"2.82206E-06"
# Row 23: l238
$ws.Range("A23").Value = "l238"
$ws.Range("B23").Value = [double]"1.55125E-10"
# Row 24: NA
$ws.Range("A24").Value = "NA"
$ws.Range("B24").Value = [double]"6.02214129E+23"
# Row 25: NR85
$ws.Range("A25").Value = "NR85"
$ws.Range("B25").Value = 137.881
# Row 26: cps
$ws.Range("A26").Value = "cps"
$ws.Range("B26").Value = 62500000
# Row 27: slope
$ws.Range("A27").Value = "slope"
$ws.Range("B27").Value = 3.4053
# Row 28: R3433
$ws.Range("A28").Value = "R3433"
$ws.Range("B28").Value = 0.002324
# Row 29: R3533
$ws.Range("A29").Value = "R3533"
$ws.Range("B29").Value = 0.005066
# Row 30: R3029
$ws.Range("A30").Value = "R3029"
$ws.Range("B30").Value = 0.00005
# Row 31: th229SubU238
$ws.Range("A31").Value = "th229SubU238"
$ws.Range("B31").Value = 0
# Row 32: th230SubU238
$ws.Range("A32").Value = "th230SubU238"
$ws.Range("B32").Value = 0
# Row 33: tri229
$ws.Range("A33").Value = "tri229"
$ws.Range("B33").Value = 0.018067
# Row 34: tri233
$ws.Range("A34").Value = "tri233"
$ws.Range("B34").Value = 0.038556
# Row 35: tri236
$ws.Range("A35").Value = "tri236"
$ws.Range("B35").Value = 3.86778
# Row 36: blank232
$ws.Range("A36").Value = "blank232"
$ws.Range("B36").Value = 0.003
# Row 37: blank232S
$ws.Range("A37").Value = "blank232S"
$ws.Range("B37").Value = 0
# Row 38: blank234
$ws.Range("A38").Value = "blank234"
$ws.Range("B38").Value = 0.3
# Row 39: blank234S
$ws.Range("A39").Value = "blank234S"
$ws.Range("B39").Value = 0
# Row 40: blank238
$ws.Range("A40").Value = "blank238"
$ws.Range("B40").Value = 0.005
# Row 41: blank238S
$ws.Range("A41").Value = "blank238S"
$ws.Range("B41").Value = 0
# Row 42: ch_blank230
$ws.Range("A42").Value = "ch_blank230"
$ws.Range("B42").Value = 0.05
# Row 43: ch_blank230S
$ws.Range("A43").Value = "ch_blank230S"
$ws.Range("B43").Value = 0
# Row 44: a230232_init
$ws.Range("A44").Value = "a230232_init"
$ws.Range("B44").Value = 0.75
# Row 45: a230232_init_err
$ws.Range("A45").Value = "a230232_init_err"
$ws.Range("B45").Value = 0.375
# Row 46: standardBezeich
$ws.Range("A46").Value = "standardBezeich"
$ws.Range("B46").Value = "Hu1-13"
# Row 47: standardEinwaage
$ws.Range("A47").Value = "standardEinwaage"
$ws.Range("B47").Value = 1.10995
# Row 48: standardTriSp13
$ws.Range("A48").Value = "standardTriSp13"
$ws.Range("B48").Value = 10.34908
# Row 49: type
$ws.Range("A49").Value = "type"
$ws.Range("B49").Value = "stalag"

# ---------------------------------------------------------------------
# 2) Calc sheet - re-computed outputs for the 4 samples (rows 7,9,11,13)
# ---------------------------------------------------------------------
$calc = $wb.Worksheets.Item("Calc")

$calc.Range("AP7").Value = 270.1318
$calc.Range("AQ7").Value = 41.96341426094132
$calc.Range("AW7").Value = 344.5467
$calc.Range("AY7").Value = 53.52348865572346
$calc.Range("BC7").Value = 12.61180836448982
$calc.Range("BG7").Value = 78.96000000000001
$calc.Range("BH7").Value = 78.34

$calc.Range("AP9").Value = 298.3195
$calc.Range("AQ9").Value = 44.06339574082669
$calc.Range("AW9").Value = 400.3249
$calc.Range("AY9").Value = 59.13030826938288
$calc.Range("BC9").Value = 10.88748143802084
$calc.Range("BG9").Value = 73.86
$calc.Range("BH9").Value = 73.72

$calc.Range("AP11").Value = 288.6347
$calc.Range("AQ11").Value = 44.08018942123226
$calc.Range("AW11").Value = 388.7997
$calc.Range("AY11").Value = 59.37752513123618
$calc.Range("BC11").Value = 15.58068770008644
$calc.Range("BG11").Value = 74.16
$calc.Range("BH11").Value = 73.7

$calc.Range("AP13").Value = 212.0804
$calc.Range("AQ13").Value = 35.37932358691081
$calc.Range("AW13").Value = 230.8497
$calc.Range("AY13").Value = 38.51054015752892
$calc.Range("BC13").Value = 8.274529688906103
$calc.Range("BG13").Value = 88.28
$calc.Range("BH13").Value = 88.86

# ---------------------------------------------------------------------
# 3) Results sheet - duplicates of AP/AW/BC from Calc (columns N/P/R)
# ---------------------------------------------------------------------
$results = $wb.Worksheets.Item("Results")

$results.Range("N7").Value = 270.1318
$results.Range("P7").Value = 344.5467
$results.Range("R7").Value = 12.61180836448982

$results.Range("N9").Value = 298.3195
$results.Range("P9").Value = 400.3249
$results.Range("R9").Value = 10.88748143802084

$results.Range("N11").Value = 288.6347
$results.Range("P11").Value = 388.7997
$results.Range("R11").Value = 15.58068770008644

$results.Range("N13").Value = 212.0804
$results.Range("P13").Value = 230.8497
$results.Range("R13").Value = 8.274529688906103
